$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 2879.4167
$ws.Range("I40").Value = 3045.3
$ws.Range("J40").Value = 2050
$ws.Range("K40").Value = 3045.3
$ws.Range("L40").Value = 2050
$ws.Range("M40").Value = -2870.3
$ws.Range("N40").Value = -2400

# Row 55 (Leve Item ID 5517)
$ws.Range("H55").Value = 73.375
$ws.Range("I55").Value = 73.375
$ws.Range("K55").Value = 73.375
$ws.Range("M55").Value = 140.625

# Row 64 (Leve Item ID 5506)
$ws.Range("H64").Value = 271916.12
$ws.Range("I64").Value = 355062.56
$ws.Range("J64").Value = 3999.889
$ws.Range("K64").Value = 355062.56
$ws.Range("L64").Value = 3999.889
$ws.Range("M64").Value = -354814.56
$ws.Range("N64").Value = -4495.889

# Row 67 (Leve Item ID 5506)
$ws.Range("H67").Value = 271916.12
$ws.Range("I67").Value = 355062.56
$ws.Range("J67").Value = 3999.889
$ws.Range("K67").Value = 355062.56
$ws.Range("L67").Value = 3999.889
$ws.Range("M67").Value = -354204.56
$ws.Range("N67").Value = -5715.889

# Row 76 (Leve Item ID 12602)
$ws.Range("H76").Value = 3317.869
$ws.Range("I76").Value = 3214.182
$ws.Range("J76").Value = 4268.3335
$ws.Range("K76").Value = 3214.182
$ws.Range("L76").Value = 4268.3335
$ws.Range("M76").Value = -2899.182
$ws.Range("N76").Value = -4898.3335

# Row 79 (Leve Item ID 12602)
$ws.Range("H79").Value = 3317.869
$ws.Range("I79").Value = 3214.182
$ws.Range("J79").Value = 4268.3335
$ws.Range("K79").Value = 3214.182
$ws.Range("L79").Value = 4268.3335
$ws.Range("M79").Value = -2122.182
$ws.Range("N79").Value = -6452.3335

# Row 80 (Leve Item ID 12605)
$ws.Range("H80").Value = 5360.0415
$ws.Range("I80").Value = 409.15384
$ws.Range("J80").Value = 11211.091
$ws.Range("K80").Value = 1227.46152
$ws.Range("L80").Value = 33633.273
$ws.Range("M80").Value = -229.4615200000001
$ws.Range("N80").Value = -35629.273

# Row 83 (Leve Item ID 12605)
$ws.Range("H83").Value = 5360.0415
$ws.Range("I83").Value = 409.15384
$ws.Range("J83").Value = 11211.091
$ws.Range("K83").Value = 3682.38456
$ws.Range("L83").Value = 100899.819
$ws.Range("M83").Value = 1309.61544
$ws.Range("N83").Value = -110883.819

# Row 135 (Leve Item ID 44047)
$ws.Range("H135").Value = 46876290
$ws.Range("I135").Value = 20001452
$ws.Range("J135").Value = 142857840
$ws.Range("K135").Value = 180013068
$ws.Range("L135").Value = 1285720560
$ws.Range("M135").Value = -180010533
$ws.Range("N135").Value = -1285725630

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2684.3125
$ws.Range("I138").Value = 1244.5227
$ws.Range("J138").Value = 3902.5962
$ws.Range("K138").Value = 3733.5681
$ws.Range("L138").Value = 11707.7886
$ws.Range("M138").Value = 1406.4319
$ws.Range("N138").Value = -21987.7886

$ws = $wb.Worksheets.Item("ARM")
# Row 63 (Leve Item ID 12528)
$ws.Range("H63").Value = 2201
$ws.Range("I63").Value = 1501.25
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 1501.25
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -815.25
$ws.Range("N63").Value = -6372

# Row 66 (Leve Item ID 12528)
$ws.Range("H66").Value = 2201
$ws.Range("I66").Value = 1501.25
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 7506.25
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -4074.25
$ws.Range("N66").Value = -31864

# Row 88 (Leve Item ID 12530)
$ws.Range("H88").Value = 5467.5713
$ws.Range("J88").Value = 2176.2856
$ws.Range("L88").Value = 2176.2856
$ws.Range("N88").Value = -2988.2856

# Row 91 (Leve Item ID 12530)
$ws.Range("H91").Value = 5467.5713
$ws.Range("J91").Value = 2176.2856
$ws.Range("L91").Value = 2176.2856
$ws.Range("N91").Value = -4984.2856

# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 2819.889
$ws.Range("I102").Value = 2736
$ws.Range("J102").Value = 2924.75
$ws.Range("K102").Value = 2736
$ws.Range("L102").Value = 2924.75
$ws.Range("M102").Value = -1114
$ws.Range("N102").Value = -6168.75

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 995.8
$ws.Range("I20").Value = 978.1111
$ws.Range("J20").Value = 1041.2858
$ws.Range("K20").Value = 978.1111
$ws.Range("L20").Value = 1041.2858
$ws.Range("M20").Value = -731.1111
$ws.Range("N20").Value = -1535.2858

# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 4625.617
$ws.Range("I105").Value = 4317.8276
$ws.Range("J105").Value = 5121.5
$ws.Range("K105").Value = 4317.8276
$ws.Range("L105").Value = 5121.5
$ws.Range("M105").Value = -2570.8276
$ws.Range("N105").Value = -8615.5

$ws = $wb.Worksheets.Item("CUL")
# Row 80 (Leve Item ID 12890)
$ws.Range("H80").Value = 4544.4443
$ws.Range("J80").Value = 4571.4287
$ws.Range("L80").Value = 13714.2861
$ws.Range("N80").Value = -15586.2861

# Row 83 (Leve Item ID 12890)
$ws.Range("H83").Value = 4544.4443
$ws.Range("J83").Value = 4571.4287
$ws.Range("L83").Value = 41142.85830000001
$ws.Range("N83").Value = -50502.85830000001

# Row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 1255.6571
$ws.Range("I107").Value = 349.77777
$ws.Range("J107").Value = 1569.2307
$ws.Range("K107").Value = 1049.33331
$ws.Range("L107").Value = 4707.6921
$ws.Range("M107").Value = 870.66669
$ws.Range("N107").Value = -8547.6921

# Row 139 (Leve Item ID 44102)
$ws.Range("H139").Value = 1762824.5
$ws.Range("I139").Value = 2710419.2
$ws.Range("K139").Value = 8131257.600000001
$ws.Range("M139").Value = -8126117.600000001

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 4569
$ws.Range("I70").Value = 3953.8572
$ws.Range("J70").Value = 5231.4614
$ws.Range("K70").Value = 3953.8572
$ws.Range("L70").Value = 5231.4614
$ws.Range("M70").Value = -3683.8572
$ws.Range("N70").Value = -5771.4614

# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 4569
$ws.Range("I73").Value = 3953.8572
$ws.Range("J73").Value = 5231.4614
$ws.Range("K73").Value = 3953.8572
$ws.Range("L73").Value = 5231.4614
$ws.Range("M73").Value = -3017.8572
$ws.Range("N73").Value = -7103.4614

# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 9543.25
$ws.Range("I80").Value = 26175
$ws.Range("J80").Value = 3999.3333
$ws.Range("K80").Value = 26175
$ws.Range("L80").Value = 3999.3333
$ws.Range("M80").Value = -25177
$ws.Range("N80").Value = -5995.3333

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 9543.25
$ws.Range("I83").Value = 26175
$ws.Range("J83").Value = 3999.3333
$ws.Range("K83").Value = 130875
$ws.Range("L83").Value = 19996.6665
$ws.Range("M83").Value = -125883
$ws.Range("N83").Value = -29980.6665

$ws = $wb.Worksheets.Item("LTW")
# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 3796.524
$ws.Range("I100").Value = 2393.3635
$ws.Range("J100").Value = 5340
$ws.Range("K100").Value = 2393.3635
$ws.Range("L100").Value = 5340
$ws.Range("M100").Value = -1852.3635
$ws.Range("N100").Value = -6422

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (Leve Item ID 12596)
$ws.Range("H81").Value = 1506.3572
$ws.Range("I81").Value = 886.125
$ws.Range("J81").Value = 2333.3333
$ws.Range("K81").Value = 1772.25
$ws.Range("L81").Value = 4666.6666
$ws.Range("M81").Value = -711.25
$ws.Range("N81").Value = -6788.6666

# Row 84 (Leve Item ID 12596)
$ws.Range("H84").Value = 1506.3572
$ws.Range("I84").Value = 886.125
$ws.Range("J84").Value = 2333.3333
$ws.Range("K84").Value = 8861.25
$ws.Range("L84").Value = 23333.333
$ws.Range("M84").Value = -3557.25
$ws.Range("N84").Value = -33941.333
